# Updated via Streamlit Approval System
# Move the pending "HOLD" flag off rows 2, 3 and 8 and onto row 7
# in the APPROVAL_1 (AI) / APPROVAL_2 (AJ) columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Clear HOLD from rows 2, 3 and 8
$ws.Range("AI2").Value = ""
$ws.Range("AJ2").Value = ""

$ws.Range("AI3").Value = ""
$ws.Range("AJ3").Value = ""

$ws.Range("AI8").Value = ""
$ws.Range("AJ8").Value = ""

# Put row 7 on HOLD
$ws.Range("AI7").Value = "HOLD"
$ws.Range("AJ7").Value = "HOLD"
